$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header row: Id, Name, Client, Type, File, Date
$ws.Range("A1").Value = "Id"
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "Client"
$ws.Range("D1").Value = "Type"
$ws.Range("E1").Value = "File"
$ws.Range("F1").Value = "Date"

$date = "Feb 6, 2022 (05:00:11 EST)"

$rows = @(
    @{ Id = 1;  Name = "Supplier";     Client = "DRX"; Type = "I"; File = "I_Supplier.xml" },
    @{ Id = 2;  Name = "Plant";        Client = "DRX"; Type = "I"; File = "I_Plant.xml" },
    @{ Id = 3;  Name = "Solicitation"; Client = "DRX"; Type = "I"; File = "I_Solicitation.xml" },
    @{ Id = 4;  Name = "BOM";          Client = "DRX"; Type = "I"; File = "I_BOM.xml" },
    @{ Id = 5;  Name = "RequestFile";  Client = "DRX"; Type = "E"; File = "E_RequestFile.xml" },
    @{ Id = 6;  Name = "BOM";          Client = "DRX"; Type = "E"; File = "E_BOM.xml" },
    @{ Id = 7;  Name = "Supplier";     Client = "GYU"; Type = "I"; File = "I_Supplier.xml" },
    @{ Id = 8;  Name = "Plant";        Client = "GYU"; Type = "I"; File = "I_Plant.xml" },
    @{ Id = 9;  Name = "Solicitation"; Client = "GYU"; Type = "I"; File = "I_Solicitation.xml" },
    @{ Id = 10; Name = "BOM";          Client = "GYU"; Type = "I"; File = "I_BOM.xml" },
    @{ Id = 11; Name = "RequestFile";  Client = "GYU"; Type = "E"; File = "E_RequestFile.xml" },
    @{ Id = 12; Name = "BOM";          Client = "GYU"; Type = "E"; File = "E_BOM.xml" }
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row.Id
    $ws.Cells.Item($r, 2).Value = $row.Name
    $ws.Cells.Item($r, 3).Value = $row.Client
    $ws.Cells.Item($r, 4).Value = $row.Type
    $ws.Cells.Item($r, 5).Value = $row.File
    $ws.Cells.Item($r, 6).Value = $date
    $r = $r + 1
}
